$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Employee_name, Netpay, Month) for rows 2-5
$data = @(
    @("Vidya Sagar  Pogiri", 26454.6, "November"),
    @("Balaraju vankala",    35856.5, "November"),
    @("Priyanka Muddana",    49460.8, "November"),
    @("pattabhi ramarao",    8701.46, "November")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    # Keep column A present-but-empty (matches the original template's blank
    # Account Number column) without clobbering any existing cell.
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
